# Fix formatting of scraped floating point numbers (and a few stray
# commas used as name separators) that were mis-rendered with
# Spanish/European-style grouping ("." thousands, "," decimal) instead of
# plain decimal notation ("." decimal, no thousands separator) when the
# source data was originally scraped into this workbook.
#
# Column H ("Importe") holds amounts such as "725,50" / "964.380,00" that
# must become "725.50" / "964380.00". A handful of "Razon social"/
# "Nombre Fantasia" entries (column E/F) also had a comma used where a
# period was intended as a separator between two co-holders' names, e.g.
# "FERNANDEZ MARIO H, GALLICET OSCAR M" -> "FERNANDEZ MARIO H. GALLICET OSCAR M".
#
# Each entry below is the target cell plus its corrected literal text. For
# values that are purely numeric-looking, the new text is prefixed with a
# leading apostrophe: this is how Excel's own UI keeps/enters a
# number-shaped value as literal text (a *quote-prefixed* text cell)
# instead of silently re-parsing "725.50" back into the number 725.5 (which
# would also lose the trailing zero and stop matching the source text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = "E67"; Value = "FERNANDEZ MARIO H. GALLICET OSCAR M" },
    @{ Cell = "E88"; Value = "FERNANDEZ MARIO H. GALLICET OSCAR M" },
    @{ Cell = "E130"; Value = "FERNANDEZ MARIO H. GALLICET OSCAR M" },
    @{ Cell = "E69"; Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA" },
    @{ Cell = "F69"; Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA" },
    @{ Cell = "E70"; Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO" },
    @{ Cell = "E95"; Value = "RICCOTTI. MARIANA EDITH" },
    @{ Cell = "E107"; Value = "RICCOTTI. MARIANA EDITH" },
    @{ Cell = "H2"; Value = "'725.50" },
    @{ Cell = "H3"; Value = "'2655.00" },
    @{ Cell = "H4"; Value = "'2880.00" },
    @{ Cell = "H60"; Value = "'2880.00" },
    @{ Cell = "H5"; Value = "'6150.00" },
    @{ Cell = "H6"; Value = "'850.40" },
    @{ Cell = "H7"; Value = "'964380.00" },
    @{ Cell = "H8"; Value = "'1895.00" },
    @{ Cell = "H9"; Value = "'144.00" },
    @{ Cell = "H10"; Value = "'15325.07" },
    @{ Cell = "H11"; Value = "'92613.19" },
    @{ Cell = "H12"; Value = "'14.00" },
    @{ Cell = "H13"; Value = "'6785.90" },
    @{ Cell = "H14"; Value = "'597.00" },
    @{ Cell = "H15"; Value = "'12580.07" },
    @{ Cell = "H16"; Value = "'223.60" },
    @{ Cell = "H17"; Value = "'1662.90" },
    @{ Cell = "H18"; Value = "'6556.75" },
    @{ Cell = "H19"; Value = "'10015.00" },
    @{ Cell = "H20"; Value = "'239.11" },
    @{ Cell = "H21"; Value = "'69.16" },
    @{ Cell = "H22"; Value = "'1113.00" },
    @{ Cell = "H23"; Value = "'1730.00" },
    @{ Cell = "H24"; Value = "'10157.60" },
    @{ Cell = "H25"; Value = "'98.50" },
    @{ Cell = "H26"; Value = "'1786.00" },
    @{ Cell = "H27"; Value = "'7873.50" },
    @{ Cell = "H28"; Value = "'962.00" },
    @{ Cell = "H29"; Value = "'1039.52" },
    @{ Cell = "H30"; Value = "'39792.22" },
    @{ Cell = "H31"; Value = "'4190.00" },
    @{ Cell = "H32"; Value = "'382.00" },
    @{ Cell = "H33"; Value = "'3565.97" },
    @{ Cell = "H34"; Value = "'712.45" },
    @{ Cell = "H35"; Value = "'1158.73" },
    @{ Cell = "H36"; Value = "'1110.00" },
    @{ Cell = "H37"; Value = "'996.90" },
    @{ Cell = "H38"; Value = "'10881.00" },
    @{ Cell = "H39"; Value = "'27.00" },
    @{ Cell = "H40"; Value = "'336.00" },
    @{ Cell = "H41"; Value = "'310.24" },
    @{ Cell = "H42"; Value = "'21.06" },
    @{ Cell = "H43"; Value = "'1918.62" },
    @{ Cell = "H44"; Value = "'1968.18" },
    @{ Cell = "H45"; Value = "'5358.50" },
    @{ Cell = "H46"; Value = "'125.71" },
    @{ Cell = "H47"; Value = "'243.00" },
    @{ Cell = "H48"; Value = "'11335.00" },
    @{ Cell = "H49"; Value = "'16.80" },
    @{ Cell = "H50"; Value = "'3304.10" },
    @{ Cell = "H51"; Value = "'174.55" },
    @{ Cell = "H52"; Value = "'120.00" },
    @{ Cell = "H53"; Value = "'1432.81" },
    @{ Cell = "H54"; Value = "'4020.00" },
    @{ Cell = "H55"; Value = "'5905.00" },
    @{ Cell = "H56"; Value = "'15216.00" },
    @{ Cell = "H57"; Value = "'1495.00" },
    @{ Cell = "H58"; Value = "'1260.00" },
    @{ Cell = "H59"; Value = "'4485.00" },
    @{ Cell = "H61"; Value = "'3545.00" },
    @{ Cell = "H62"; Value = "'400.00" },
    @{ Cell = "H63"; Value = "'53748.00" },
    @{ Cell = "H64"; Value = "'7151.00" },
    @{ Cell = "H65"; Value = "'4000.00" },
    @{ Cell = "H66"; Value = "'159.00" },
    @{ Cell = "H67"; Value = "'2518.00" },
    @{ Cell = "H68"; Value = "'3.00" },
    @{ Cell = "H69"; Value = "'2393.15" },
    @{ Cell = "H70"; Value = "'2102.00" },
    @{ Cell = "H71"; Value = "'584.60" },
    @{ Cell = "H72"; Value = "'7901.53" },
    @{ Cell = "H73"; Value = "'2499.00" },
    @{ Cell = "H74"; Value = "'0.08" },
    @{ Cell = "H75"; Value = "'0.06" },
    @{ Cell = "H76"; Value = "'1500000.00" },
    @{ Cell = "H77"; Value = "'0.94" },
    @{ Cell = "H78"; Value = "'103769.60" },
    @{ Cell = "H79"; Value = "'1030560.00" },
    @{ Cell = "H80"; Value = "'0.70" },
    @{ Cell = "H81"; Value = "'0.02" },
    @{ Cell = "H82"; Value = "'0.82" },
    @{ Cell = "H83"; Value = "'43.80" },
    @{ Cell = "H84"; Value = "'142.39" },
    @{ Cell = "H85"; Value = "'5499.18" },
    @{ Cell = "H86"; Value = "'6594.00" },
    @{ Cell = "H87"; Value = "'550.00" },
    @{ Cell = "H88"; Value = "'9080.10" },
    @{ Cell = "H89"; Value = "'850.00" },
    @{ Cell = "H90"; Value = "'4161.79" },
    @{ Cell = "H91"; Value = "'175.70" },
    @{ Cell = "H92"; Value = "'3090.00" },
    @{ Cell = "H93"; Value = "'747.00" },
    @{ Cell = "H94"; Value = "'1327.00" },
    @{ Cell = "H95"; Value = "'12000.00" },
    @{ Cell = "H156"; Value = "'12000.00" },
    @{ Cell = "H96"; Value = "'570.00" },
    @{ Cell = "H97"; Value = "'1726.00" },
    @{ Cell = "H98"; Value = "'2500.00" },
    @{ Cell = "H99"; Value = "'55000.00" },
    @{ Cell = "H100"; Value = "'12900.00" },
    @{ Cell = "H101"; Value = "'1000.00" },
    @{ Cell = "H150"; Value = "'1000.00" },
    @{ Cell = "H102"; Value = "'60338.70" },
    @{ Cell = "H103"; Value = "'6900.00" },
    @{ Cell = "H104"; Value = "'950.00" },
    @{ Cell = "H105"; Value = "'8350.00" },
    @{ Cell = "H106"; Value = "'60750.00" },
    @{ Cell = "H107"; Value = "'10000.00" },
    @{ Cell = "H108"; Value = "'108.82" },
    @{ Cell = "H109"; Value = "'28.38" },
    @{ Cell = "H110"; Value = "'510.20" },
    @{ Cell = "H111"; Value = "'108.00" },
    @{ Cell = "H112"; Value = "'51233.50" },
    @{ Cell = "H113"; Value = "'560.00" },
    @{ Cell = "H114"; Value = "'1370.00" },
    @{ Cell = "H115"; Value = "'250.00" },
    @{ Cell = "H116"; Value = "'6023.87" },
    @{ Cell = "H117"; Value = "'6050.00" },
    @{ Cell = "H118"; Value = "'750.00" },
    @{ Cell = "H126"; Value = "'750.00" },
    @{ Cell = "H119"; Value = "'450.00" },
    @{ Cell = "H120"; Value = "'2000.00" },
    @{ Cell = "H121"; Value = "'15088.80" },
    @{ Cell = "H122"; Value = "'1400.00" },
    @{ Cell = "H123"; Value = "'1200.00" },
    @{ Cell = "H124"; Value = "'1750.00" },
    @{ Cell = "H125"; Value = "'1105.00" },
    @{ Cell = "H127"; Value = "'8720.00" },
    @{ Cell = "H128"; Value = "'436.00" },
    @{ Cell = "H129"; Value = "'873.00" },
    @{ Cell = "H130"; Value = "'1380.00" },
    @{ Cell = "H131"; Value = "'384.00" },
    @{ Cell = "H132"; Value = "'49.20" },
    @{ Cell = "H133"; Value = "'1455.00" },
    @{ Cell = "H134"; Value = "'1847.39" },
    @{ Cell = "H135"; Value = "'2326.00" },
    @{ Cell = "H136"; Value = "'132.00" },
    @{ Cell = "H137"; Value = "'733.00" },
    @{ Cell = "H138"; Value = "'7.50" },
    @{ Cell = "H139"; Value = "'980.00" },
    @{ Cell = "H140"; Value = "'23501.49" },
    @{ Cell = "H141"; Value = "'1049.00" },
    @{ Cell = "H142"; Value = "'660.18" },
    @{ Cell = "H143"; Value = "'568.75" },
    @{ Cell = "H144"; Value = "'367.00" },
    @{ Cell = "H145"; Value = "'128.70" },
    @{ Cell = "H146"; Value = "'12338.30" },
    @{ Cell = "H147"; Value = "'982.99" },
    @{ Cell = "H148"; Value = "'2659.60" },
    @{ Cell = "H149"; Value = "'2305.66" },
    @{ Cell = "H151"; Value = "'7000.00" },
    @{ Cell = "H152"; Value = "'690.00" },
    @{ Cell = "H153"; Value = "'190000.00" },
    @{ Cell = "H164"; Value = "'190000.00" },
    @{ Cell = "H154"; Value = "'69500.00" },
    @{ Cell = "H155"; Value = "'69500.00" },
    @{ Cell = "H160"; Value = "'69500.00" },
    @{ Cell = "H157"; Value = "'17000.00" },
    @{ Cell = "H158"; Value = "'137732.00" },
    @{ Cell = "H159"; Value = "'106460.00" },
    @{ Cell = "H161"; Value = "'124940.00" },
    @{ Cell = "H162"; Value = "'232386.00" },
    @{ Cell = "H163"; Value = "'134208.00" },
    @{ Cell = "H165"; Value = "'180000.00" },
    @{ Cell = "H166"; Value = "'98012.00" },
    @{ Cell = "H167"; Value = "'229295.87" },
    @{ Cell = "H168"; Value = "'69400.00" },
    @{ Cell = "H169"; Value = "'1925.00" },
    @{ Cell = "H170"; Value = "'1040.85" },
    @{ Cell = "H171"; Value = "'1700.00" }

)

foreach ($edit in $edits) {
    $ws.Range($edit.Cell).Value = $edit.Value
}
